$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of daily-expense data (row 21), copying the date
# cell's number formatting from the row above (row 20) so it reuses the
# existing style rather than minting a new one.
$ws.Range("A20").Copy()
$ws.Range("A21").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A21").Value = 43809

$ws.Range("B21").Value = 0
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 12.5
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 3

$ws.Range("K21").Select()
